$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 479, shifting every
# subsequent row (old 479..588) down to (481..590). This reproduces the
# observed diff: every row from new-480.. onward equals the old row two
# positions above it, and the two rows that fall off the bottom of the
# original range reappear (duplicated) as the new last two rows 589/590.
$ws.Range("A479:A480").EntireRow.Insert()

# New row 479: weekly "Primera" quality entry for the new date.
$ws.Range("A479").Value = 1
$ws.Range("B479").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C479").Value = "Arica y Parinacota"
$ws.Range("D479").Value = 45173
$ws.Range("E479").Value = 15
$ws.Range("F479").Value = 100112032
$ws.Range("G479").Value = "Zapallo italiano"
$ws.Range("H479").Value = "Huracán"
$ws.Range("I479").Value = "Primera"
$ws.Range("J479").Value = 130
$ws.Range("K479").Value = 7000
$ws.Range("L479").Value = 8000
$ws.Range("M479").Value = 7500
$ws.Range("N479").Value = "$/caja 70 unidades"
$ws.Range("O479").Value = "Región de Arica y Parinacota"
$ws.Range("P479").Value = 107
$ws.Range("Q479").Value = 70
$ws.Range("R479").Value = "Hortaliza"

# New row 480: weekly "Segunda" quality entry for the same new date.
$ws.Range("A480").Value = 1
$ws.Range("B480").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C480").Value = "Arica y Parinacota"
$ws.Range("D480").Value = 45173
$ws.Range("E480").Value = 15
$ws.Range("F480").Value = 100112032
$ws.Range("G480").Value = "Zapallo italiano"
$ws.Range("H480").Value = "Huracán"
$ws.Range("I480").Value = "Segunda"
$ws.Range("J480").Value = 150
$ws.Range("K480").Value = 5000
$ws.Range("L480").Value = 6000
$ws.Range("M480").Value = 5500
$ws.Range("N480").Value = "$/caja 100 unidades"
$ws.Range("O480").Value = "Región de Arica y Parinacota"
$ws.Range("P480").Value = 55
$ws.Range("Q480").Value = 100
$ws.Range("R480").Value = "Hortaliza"
